$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch to manual calculation so that we can poke the cached results of
# the "What-If" data table formulas (G2:G5 / J2:J5) directly - those are
# driven by an Excel Data Table feature that this runtime does not
# actually re-simulate, so the cached <v> has to be supplied explicitly.
$excel.Calculation = -4135   # xlCalculationManual

# --- Row 2 (NE @ KC) ---
$ws.Range("B2").Value = "KC"
$ws.Range("C2").Value = -9.5
$ws.Range("D2").Value = 37.5

# --- Row 3 (KC @ KC) ---
$ws.Range("C3").Value = -9.5
$ws.Range("D3").Value = -37.5

# --- Row 4 (MIN @ MIN) ---
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()

# --- Row 5 (LAR @ LAR) ---
$ws.Range("B5").Value = "LAR"
$ws.Range("D5").Value = 49

# Recalculate the ordinary formulas (E,F,H,I columns, B9:C12 strings, etc.)
# while leaving the manual data-table cache values untouched.
$ws.Calculate()

# Write the recalculated "What-If" data table cached results (G2:G5 /
# J2:J5) that the real Excel Data Table feature would have produced.
$ws.Cells.Item(2, 7).Value = 24
$ws.Cells.Item(3, 7).Value = 17
$ws.Cells.Item(4, 7).Value = 17
$ws.Cells.Item(5, 7).Value = 28

$ws.Cells.Item(2, 10).Value = 14
$ws.Cells.Item(3, 10).Value = 3
$ws.Cells.Item(4, 10).Value = 3
$ws.Cells.Item(5, 10).Value = 21

# Recalculate once more so the strings that reference G/J (B9:B12, C10:C12)
# pick up the new scores.
$ws.Calculate()

# Match the saved selection (the author had C4 selected when they saved).
$ws.Range("C4").Select()
